# Scheduled market-data refresh: updates Universalis price/profit
# columns (H:N) for the rows whose current-average-price snapshot changed.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 81.07143000000001
$ws.Range("I55").Value = 57.5
$ws.Range("J55").Value = 85
$ws.Range("K55").Value = 57.5
$ws.Range("L55").Value = 85
$ws.Range("M55").Value = 156.5
$ws.Range("N55").Value = -513
$ws.Range("H86").Value = 5637.727
$ws.Range("J86").Value = 15248.857
$ws.Range("L86").Value = 15248.857
$ws.Range("N86").Value = -17494.857
$ws.Range("H89").Value = 5637.727
$ws.Range("J89").Value = 15248.857
$ws.Range("L89").Value = 76244.285
$ws.Range("N89").Value = -87476.285
$ws.Range("H98").Value = 695.8261
$ws.Range("I98").Value = 448.26666
$ws.Range("J98").Value = 1160
$ws.Range("K98").Value = 448.26666
$ws.Range("L98").Value = 1160
$ws.Range("M98").Value = 1049.73334
$ws.Range("N98").Value = -4156
$ws.Range("H107").Value = 792.4815
$ws.Range("I107").Value = 587
$ws.Range("K107").Value = 587
$ws.Range("M107").Value = 1333
$ws.Range("H112").Value = 1019.4167
$ws.Range("I112").Value = 635.7143
$ws.Range("J112").Value = 1112.0344
$ws.Range("K112").Value = 1907.1429
$ws.Range("L112").Value = 3336.1032
$ws.Range("M112").Value = -799.1428999999998
$ws.Range("N112").Value = -5552.1032
$ws.Range("H122").Value = 695.8261
$ws.Range("I122").Value = 448.26666
$ws.Range("J122").Value = 1160
$ws.Range("K122").Value = 1344.79998
$ws.Range("L122").Value = 3480
$ws.Range("M122").Value = 1105.20002
$ws.Range("N122").Value = -8380
$ws.Range("H137").Value = 1629.0667
$ws.Range("I137").Value = 1588.6
$ws.Range("J137").Value = 1710
$ws.Range("K137").Value = 4765.799999999999
$ws.Range("L137").Value = 5130
$ws.Range("M137").Value = -2215.799999999999
$ws.Range("N137").Value = -10230
$ws.Range("H141").Value = 3404.625
$ws.Range("I141").Value = 2934
$ws.Range("K141").Value = 8802
$ws.Range("M141").Value = -3622

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5372.2925
$ws.Range("I32").Value = 4477.712
$ws.Range("K32").Value = 4477.712
$ws.Range("M32").Value = -4190.712
$ws.Range("H122").Value = 1235.2979
$ws.Range("I122").Value = 1061.5
$ws.Range("K122").Value = 3184.5
$ws.Range("M122").Value = -734.5
$ws.Range("H132").Value = 15163.368
$ws.Range("I132").Value = 1626.5333
$ws.Range("J132").Value = 65926.5
$ws.Range("K132").Value = 4879.5999
$ws.Range("L132").Value = 197779.5
$ws.Range("M132").Value = -2349.5999
$ws.Range("N132").Value = -202839.5
$ws.Range("H140").Value = 47199.6
$ws.Range("J140").Value = 47199.6
$ws.Range("L140").Value = 47199.6
$ws.Range("N140").Value = -57559.6

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3253.8865
$ws.Range("I31").Value = 2673.5
$ws.Range("J31").Value = 3585.5356
$ws.Range("K31").Value = 2673.5
$ws.Range("L31").Value = 3585.5356
$ws.Range("M31").Value = -2378.5
$ws.Range("N31").Value = -4175.5356
$ws.Range("H34").Value = 3253.8865
$ws.Range("I34").Value = 2673.5
$ws.Range("J34").Value = 3585.5356
$ws.Range("K34").Value = 2673.5
$ws.Range("L34").Value = 3585.5356
$ws.Range("M34").Value = -2471.5
$ws.Range("N34").Value = -3989.5356

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 100789.14
$ws.Range("I131").Value = 840
$ws.Range("J131").Value = 102828.92
$ws.Range("K131").Value = 2520
$ws.Range("L131").Value = 308486.76
$ws.Range("M131").Value = 2520
$ws.Range("N131").Value = -318566.76

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4176879.8
$ws.Range("I70").Value = 4123.5
$ws.Range("K70").Value = 4123.5
$ws.Range("M70").Value = -3853.5
$ws.Range("H73").Value = 4176879.8
$ws.Range("I73").Value = 4123.5
$ws.Range("K73").Value = 4123.5
$ws.Range("M73").Value = -3187.5
$ws.Range("H113").Value = 14244.429
$ws.Range("J113").Value = 6000
$ws.Range("L113").Value = 6000
$ws.Range("N113").Value = -10340
$ws.Range("H126").Value = 5665.778
$ws.Range("I126").Value = 4553.3335
$ws.Range("J126").Value = 7056.3335
$ws.Range("K126").Value = 13660.0005
$ws.Range("L126").Value = 21169.0005
$ws.Range("M126").Value = -11190.0005
$ws.Range("N126").Value = -26109.0005
$ws.Range("H132").Value = 21324.154
$ws.Range("I132").Value = 1874
$ws.Range("J132").Value = 74117.42999999999
$ws.Range("K132").Value = 5622
$ws.Range("L132").Value = 222352.29
$ws.Range("M132").Value = -3092
$ws.Range("N132").Value = -227412.29

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4223
$ws.Range("I7").Value = 4150
$ws.Range("K7").Value = 4150
$ws.Range("M7").Value = -4038
$ws.Range("H42").Value = 4083.3333
$ws.Range("J42").Value = 4083.3333
$ws.Range("L42").Value = 4083.3333
$ws.Range("N42").Value = -5209.3333
$ws.Range("H49").Value = 4083.3333
$ws.Range("J49").Value = 4083.3333
$ws.Range("L49").Value = 4083.3333
$ws.Range("N49").Value = -4377.3333
$ws.Range("H68").Value = 1968.4667
$ws.Range("I68").Value = 1547.2858
$ws.Range("J68").Value = 2337
$ws.Range("K68").Value = 1547.2858
$ws.Range("L68").Value = 2337
$ws.Range("M68").Value = -798.2858000000001
$ws.Range("N68").Value = -3835
$ws.Range("H71").Value = 1968.4667
$ws.Range("I71").Value = 1547.2858
$ws.Range("J71").Value = 2337
$ws.Range("K71").Value = 7736.429
$ws.Range("L71").Value = 11685
$ws.Range("M71").Value = -3992.429
$ws.Range("N71").Value = -19173
$ws.Range("H93").Value = 1051.5
$ws.Range("I93").Value = 801
$ws.Range("K93").Value = 801
$ws.Range("M93").Value = 447
$ws.Range("H122").Value = 579199.7
$ws.Range("I122").Value = 1091395.2
$ws.Range("J122").Value = 2979.6875
$ws.Range("K122").Value = 3274185.6
$ws.Range("L122").Value = 8939.0625
$ws.Range("M122").Value = -3271735.6
$ws.Range("N122").Value = -13839.0625
$ws.Range("H126").Value = 4223
$ws.Range("I126").Value = 4150
$ws.Range("K126").Value = 12450
$ws.Range("M126").Value = -9980
$ws.Range("H127").Value = 37442.79
$ws.Range("J127").Value = 37442.79
$ws.Range("L127").Value = 37442.79
$ws.Range("N127").Value = -47362.79
$ws.Range("H136").Value = 1987.75
$ws.Range("I136").Value = 2100.2856
$ws.Range("J136").Value = 1200
$ws.Range("K136").Value = 6300.8568
$ws.Range("L136").Value = 3600
$ws.Range("M136").Value = -3750.8568
$ws.Range("N136").Value = -8700

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1227
$ws.Range("I126").Value = 1255.5333
$ws.Range("J126").Value = 799
$ws.Range("K126").Value = 3766.5999
$ws.Range("L126").Value = 2397
$ws.Range("M126").Value = -1296.5999
$ws.Range("N126").Value = -7337
$ws.Range("H136").Value = 25809126
$ws.Range("I136").Value = 41291670
$ws.Range("K136").Value = 123875010
$ws.Range("M136").Value = -123872460
